$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.03478265776327
$ws.Cells.Item(2, 4).Value = 1.056439223162817
$ws.Cells.Item(2, 5).Value = 1.045639981478569
$ws.Cells.Item(2, 6).Value = 1.060715653429755
$ws.Cells.Item(2, 9).Value = 1.04399115998102
$ws.Cells.Item(2, 10).Value = 1.03990001639786
$ws.Cells.Item(2, 11).Value = 1.059176931157509
$ws.Cells.Item(2, 12).Value = 1.048407655049973
$ws.Cells.Item(2, 13).Value = 1.063441677879617
$ws.Cells.Item(2, 14).Value = 1.017259493698236
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.035558520789289
$ws.Cells.Item(3, 4).Value = 1.057039375262869
$ws.Cells.Item(3, 5).Value = 1.046311645168585
$ws.Cells.Item(3, 6).Value = 1.061426131592645
$ws.Cells.Item(3, 9).Value = 1.044158324048745
$ws.Cells.Item(3, 10).Value = 1.04032003246765
$ws.Cells.Item(3, 11).Value = 1.059591266142113
$ws.Cells.Item(3, 12).Value = 1.04889116404999
$ws.Cells.Item(3, 13).Value = 1.063966898297183
$ws.Cells.Item(3, 14).Value = 1.017400035706474
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.036061400031598
$ws.Cells.Item(4, 4).Value = 1.057428314850099
$ws.Cells.Item(4, 5).Value = 1.046747347663646
$ws.Cells.Item(4, 6).Value = 1.061886896247059
$ws.Cells.Item(4, 9).Value = 1.044265713028606
$ws.Cells.Item(4, 10).Value = 1.040591963183721
$ws.Cells.Item(4, 11).Value = 1.059859270603514
$ws.Cells.Item(4, 12).Value = 1.049204433836831
$ws.Cells.Item(4, 13).Value = 1.064307111419134
$ws.Cells.Item(4, 14).Value = 1.017490992804391
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.036273010562365
$ws.Cells.Item(5, 4).Value = 1.057591966819859
$ws.Cells.Item(5, 5).Value = 1.046930776142306
$ws.Cells.Item(5, 6).Value = 1.062080847906192
$ws.Cells.Item(5, 9).Value = 1.044310672421256
$ws.Cells.Item(5, 10).Value = 1.040706317938769
$ws.Cells.Item(5, 11).Value = 1.059971914630743
$ws.Cells.Item(5, 12).Value = 1.049336228344248
$ws.Cells.Item(5, 13).Value = 1.064450221503958
$ws.Cells.Item(5, 14).Value = 1.017529234747424
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.036308552578637
$ws.Cells.Item(6, 4).Value = 1.057619452957134
$ws.Cells.Item(6, 5).Value = 1.04696158972713
$ws.Cells.Item(6, 6).Value = 1.062113427614989
$ws.Cells.Item(6, 9).Value = 1.044318210314922
$ws.Cells.Item(6, 10).Value = 1.040725520622915
$ws.Cells.Item(6, 11).Value = 1.059990826543411
$ws.Cells.Item(6, 12).Value = 1.049358362798995
$ws.Cells.Item(6, 13).Value = 1.064474255219268
$ws.Cells.Item(6, 14).Value = 1.017535655933386
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.036064226798068
$ws.Cells.Item(7, 4).Value = 1.057430501021356
$ws.Cells.Item(7, 5).Value = 1.046749797628293
$ws.Cells.Item(7, 6).Value = 1.061889486874058
$ws.Cells.Item(7, 9).Value = 1.044266314513778
$ws.Cells.Item(7, 10).Value = 1.040593491061417
$ws.Cells.Item(7, 11).Value = 1.059860775858943
$ws.Cells.Item(7, 12).Value = 1.049206194506697
$ws.Cells.Item(7, 13).Value = 1.064309023333361
$ws.Cells.Item(7, 14).Value = 1.017491503781665
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.035044688377481
$ws.Cells.Item(8, 4).Value = 1.056641921647652
$ws.Cells.Item(8, 5).Value = 1.045866746182823
$ws.Cells.Item(8, 6).Value = 1.060955546623104
$ws.Cells.Item(8, 9).Value = 1.044047814313205
$ws.Cells.Item(8, 10).Value = 1.040041930454942
$ws.Cells.Item(8, 11).Value = 1.059316976936046
$ws.Cells.Item(8, 12).Value = 1.048570973931707
$ws.Cells.Item(8, 13).Value = 1.06361910288165
$ws.Cells.Item(8, 14).Value = 1.017306986705752
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.033254681891243
$ws.Cells.Item(9, 4).Value = 1.055257041265797
$ws.Cells.Item(9, 5).Value = 1.044319145763123
$ws.Cells.Item(9, 6).Value = 1.059317872053371
$ws.Cells.Item(9, 9).Value = 1.043656872489259
$ws.Cells.Item(9, 10).Value = 1.039071236424284
$ws.Cells.Item(9, 11).Value = 1.058358053761258
$ws.Cells.Item(9, 12).Value = 1.047454822118665
$ws.Cells.Item(9, 13).Value = 1.062406214236996
$ws.Cells.Item(9, 14).Value = 1.016981996322435
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.032065865491123
$ws.Cells.Item(10, 4).Value = 1.054337073046282
$ws.Cells.Item(10, 5).Value = 1.043293217352586
$ws.Cells.Item(10, 6).Value = 1.058231634260733
$ws.Cells.Item(10, 9).Value = 1.043392316793262
$ws.Cells.Item(10, 10).Value = 1.03842501313363
$ws.Cells.Item(10, 11).Value = 1.057718409814591
$ws.Cells.Item(10, 12).Value = 1.046712958911728
$ws.Cells.Item(10, 13).Value = 1.061599638035473
$ws.Cells.Item(10, 14).Value = 1.016765468113919
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.031552192066767
$ws.Cells.Item(11, 4).Value = 1.05393952292968
$ws.Cells.Item(11, 5).Value = 1.042850382404515
$ws.Cells.Item(11, 6).Value = 1.057762625164856
$ws.Cells.Item(11, 9).Value = 1.043276840268712
$ws.Cells.Item(11, 10).Value = 1.038145423535725
$ws.Cells.Item(11, 11).Value = 1.057441371352355
$ws.Cells.Item(11, 12).Value = 1.04639227485631
$ws.Cells.Item(11, 13).Value = 1.061250882629225
$ws.Cells.Item(11, 14).Value = 1.016671746620976
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.031361556952013
$ws.Cells.Item(12, 4).Value = 1.053791977844804
$ws.Cells.Item(12, 5).Value = 1.042686106111955
$ws.Cells.Item(12, 6).Value = 1.057588617819391
$ws.Cells.Item(12, 9).Value = 1.043233809454741
$ws.Cells.Item(12, 10).Value = 1.038041607374687
$ws.Cells.Item(12, 11).Value = 1.057338458299776
$ws.Cells.Item(12, 12).Value = 1.046273242428553
$ws.Cells.Item(12, 13).Value = 1.061121415929637
$ws.Cells.Item(12, 14).Value = 1.016636940302935
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.031402441283992
$ws.Cells.Item(13, 4).Value = 1.053823621177185
$ws.Cells.Item(13, 5).Value = 1.042721334291482
$ws.Cells.Item(13, 6).Value = 1.057625933738641
$ws.Cells.Item(13, 9).Value = 1.043243045927514
$ws.Cells.Item(13, 10).Value = 1.038063874648126
$ws.Cells.Item(13, 11).Value = 1.057360533863142
$ws.Cells.Item(13, 12).Value = 1.04629877146506
$ws.Cells.Item(13, 13).Value = 1.061149183480702
$ws.Cells.Item(13, 14).Value = 1.01664440609581
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.031536430709335
$ws.Cells.Item(14, 4).Value = 1.053927324290539
$ws.Cells.Item(14, 5).Value = 1.042836798931055
$ws.Cells.Item(14, 6).Value = 1.057748237498226
$ws.Cells.Item(14, 9).Value = 1.043273286132125
$ws.Cells.Item(14, 10).Value = 1.038136841319713
$ws.Cells.Item(14, 11).Value = 1.057432864693917
$ws.Cells.Item(14, 12).Value = 1.046382433877525
$ws.Cells.Item(14, 13).Value = 1.06124017929838
$ws.Cells.Item(14, 14).Value = 1.016668869394163
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.031619008082278
$ws.Cells.Item(15, 4).Value = 1.05399123553013
$ws.Cells.Item(15, 5).Value = 1.042907968706097
$ws.Cells.Item(15, 6).Value = 1.057823619918734
$ws.Cells.Item(15, 9).Value = 1.043291899900367
$ws.Cells.Item(15, 10).Value = 1.038181803288724
$ws.Cells.Item(15, 11).Value = 1.057477429011265
$ws.Cells.Item(15, 12).Value = 1.046433992218718
$ws.Cells.Item(15, 13).Value = 1.061296255028842
$ws.Cells.Item(15, 14).Value = 1.016683942852307
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.03209997948226
$ws.Cells.Item(16, 4).Value = 1.054363474228165
$ws.Cells.Item(16, 5).Value = 1.043322636506314
$ws.Cells.Item(16, 6).Value = 1.058262789261596
$ws.Cells.Item(16, 9).Value = 1.043399961233795
$ws.Cells.Item(16, 10).Value = 1.038443573526719
$ws.Cells.Item(16, 11).Value = 1.057736794644069
$ws.Cells.Item(16, 12).Value = 1.046734253333594
$ws.Cells.Item(16, 13).Value = 1.061622794421574
$ws.Cells.Item(16, 14).Value = 1.016771688916563
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.032401973896437
$ws.Cells.Item(17, 4).Value = 1.05459718605967
$ws.Cells.Item(17, 5).Value = 1.04358312269776
$ws.Cells.Item(17, 6).Value = 1.058538628739017
$ws.Cells.Item(17, 9).Value = 1.043467498966268
$ws.Cells.Item(17, 10).Value = 1.03860783747844
$ws.Cells.Item(17, 11).Value = 1.057899470699533
$ws.Cells.Item(17, 12).Value = 1.04692274696471
$ws.Cells.Item(17, 13).Value = 1.061827758437845
$ws.Cells.Item(17, 14).Value = 1.016826739877772
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.032578227348118
$ws.Cells.Item(18, 4).Value = 1.054733583469024
$ws.Cells.Item(18, 5).Value = 1.043735194707252
$ws.Cells.Item(18, 6).Value = 1.058699650280215
$ws.Cells.Item(18, 9).Value = 1.043506803550386
$ws.Cells.Item(18, 10).Value = 1.038703671897129
$ws.Cells.Item(18, 11).Value = 1.057994350179366
$ws.Cells.Item(18, 12).Value = 1.047032744837317
$ws.Cells.Item(18, 13).Value = 1.061947358349471
$ws.Cells.Item(18, 14).Value = 1.016858853674032
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.032638343006001
$ws.Cells.Item(19, 4).Value = 1.054780104496823
$ws.Cells.Item(19, 5).Value = 1.043787070134611
$ws.Cells.Item(19, 6).Value = 1.058754576280906
$ws.Cells.Item(19, 9).Value = 1.043520190270432
$ws.Cells.Item(19, 10).Value = 1.038736352651984
$ws.Cells.Item(19, 11).Value = 1.058026700441679
$ws.Cells.Item(19, 12).Value = 1.0470702601628
$ws.Cells.Item(19, 13).Value = 1.061988146891446
$ws.Cells.Item(19, 14).Value = 1.016869804227098
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.032369561855987
$ws.Cells.Item(20, 4).Value = 1.054572102985503
$ws.Cells.Item(20, 5).Value = 1.043555161033836
$ws.Cells.Item(20, 6).Value = 1.058509020409439
$ws.Cells.Item(20, 9).Value = 1.043460262009657
$ws.Cells.Item(20, 10).Value = 1.038590211221699
$ws.Cells.Item(20, 11).Value = 1.057882017780629
$ws.Cells.Item(20, 12).Value = 1.046902517916415
$ws.Cells.Item(20, 13).Value = 1.061805762774002
$ws.Cells.Item(20, 14).Value = 1.016820833064721
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.031496969562715
$ws.Cells.Item(21, 4).Value = 1.053896782902742
$ws.Cells.Item(21, 5).Value = 1.042802791594007
$ws.Cells.Item(21, 6).Value = 1.057712216452574
$ws.Cells.Item(21, 9).Value = 1.043264384939676
$ws.Cells.Item(21, 10).Value = 1.038115353452251
$ws.Cells.Item(21, 11).Value = 1.057411565292148
$ws.Cells.Item(21, 12).Value = 1.046357795052854
$ws.Cells.Item(21, 13).Value = 1.061213381169025
$ws.Cells.Item(21, 14).Value = 1.016661665392319
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.030949297963578
$ws.Cells.Item(22, 4).Value = 1.053472893527528
$ws.Cells.Item(22, 5).Value = 1.042330976578848
$ws.Cells.Item(22, 6).Value = 1.057212412812605
$ws.Cells.Item(22, 9).Value = 1.043140433036031
$ws.Cells.Item(22, 10).Value = 1.037816999721388
$ws.Cells.Item(22, 11).Value = 1.057115724008204
$ws.Cells.Item(22, 12).Value = 1.046015792448184
$ws.Cells.Item(22, 13).Value = 1.060841371366029
$ws.Cells.Item(22, 14).Value = 1.016561625426567
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.031239537198069
$ws.Cells.Item(23, 4).Value = 1.053697537159065
$ws.Cells.Item(23, 5).Value = 1.042580977379252
$ws.Cells.Item(23, 6).Value = 1.057477255677552
$ws.Cells.Item(23, 9).Value = 1.043206217472195
$ws.Cells.Item(23, 10).Value = 1.037975142475043
$ws.Cells.Item(23, 11).Value = 1.057272559178514
$ws.Cells.Item(23, 12).Value = 1.046197047860106
$ws.Cells.Item(23, 13).Value = 1.061038538093791
$ws.Cells.Item(23, 14).Value = 1.01661465501078
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.032384207137068
$ws.Cells.Item(24, 4).Value = 1.054583436708638
$ws.Cells.Item(24, 5).Value = 1.043567795289892
$ws.Cells.Item(24, 6).Value = 1.058522398740707
$ws.Cells.Item(24, 9).Value = 1.04346353235405
$ws.Cells.Item(24, 10).Value = 1.038598175700514
$ws.Cells.Item(24, 11).Value = 1.057889904024389
$ws.Cells.Item(24, 12).Value = 1.046911658389897
$ws.Cells.Item(24, 13).Value = 1.061815701520047
$ws.Cells.Item(24, 14).Value = 1.016823502088532
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.033716653152181
$ws.Cells.Item(25, 4).Value = 1.055614496271106
$ws.Cells.Item(25, 5).Value = 1.044718223466604
$ws.Cells.Item(25, 6).Value = 1.059740282637269
$ws.Cells.Item(25, 9).Value = 1.04375863630167
$ws.Cells.Item(25, 10).Value = 1.039322030372275
$ws.Cells.Item(25, 11).Value = 1.058606029138482
$ws.Cells.Item(25, 12).Value = 1.047742986571177
$ws.Cells.Item(25, 13).Value = 1.062719427493189
$ws.Cells.Item(25, 14).Value = 1.017065993021982
